{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Applies the two textual corrections from the diff:\n//   1. \"wil\"   -> \"will\"   (Context Diagram paragraph)\n//   2. \"never\" -> \"newer\"  (component diagram paragraph)\n//\n// (The diff's other hunks \u2014 new w:proofErr spell-check wrappers / run\n// splits around \"Groupchats\", the lone article \"a\", and \"effect\", plus the\n// regenerated o:OLEObject ObjectID \u2014 are artifacts Word regenerates on\n// save and contain no actual text change, so there is nothing further to\n// apply here.)\n\nconst body = context.document.body;\n\n// 1) \"For my application there wil be users...\" -> \"...there will be users...\"\nlet results = body.search(\"wil be users that use the system\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"will be users that use the system\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) \"And a never version can be seen here.\" -> \"And a newer version can be seen here.\"\nresults = body.search(\"a never version can be seen here\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"a newer version can be seen here\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Applies the two textual corrections from the diff:\n#   1. \"wil\"   -> \"will\"   (Context Diagram paragraph)\n#   2. \"never\" -> \"newer\"  (component diagram paragraph)\n#\n# (The diff's other hunks -- new w:proofErr spell-check wrappers / run\n# splits around \"Groupchats\", the lone article \"a\", and \"effect\", plus the\n# regenerated o:OLEObject ObjectID -- are artifacts Word regenerates on\n# save and contain no actual text change, so there is nothing further to\n# apply here.)\n\n$d = $word.ActiveDocument\n\n# Word find/replace constants\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# 1) \"For my application there wil be users...\" -> \"...there will be users...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"wil be users that use the system\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"will be users that use the system\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n\n# 2) \"And a never version can be seen here.\" -> \"And a newer version can be seen here.\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"a never version can be seen here\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"a newer version can be seen here\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
